$wb = $excel.ActiveWorkbook

# --- Sheet: Merchant(P)CustomFee_QPayNow --- (select first so it is NOT
# left as the final active sheet/tab)
$ws2 = $wb.Worksheets.Item("Merchant(P)CustomFee_QPayNow")
$ws2.Activate()
$ws2.Range("B24").Select()

# --- Sheet: RunManager ---
$ws3 = $wb.Worksheets.Item("RunManager")
$ws3.Activate()

# Rename existing test entries (columns A) to their new names
$ws3.Range("A2").Value = "validateDashboardOnSystemLogin"
$ws3.Range("A3").Value = "validateTransactionPageOnSystemLogin"
$ws3.Range("A4").Value = "validateDashboardOnISOLogin"

# Add a new row (5) for the new ISO transaction page test, copying the
# formatting from row 4 (its closest sibling) first.
$ws3.Range("A4:D4").Copy()
$ws3.Range("A5:D5").PasteSpecial(-4122) # xlPasteFormats

$ws3.Range("A5").Value = "validateTransactionPageOnIsoLogin"
$ws3.Range("B5").Value = "To Check Whether"
$ws3.Range("C5").Value = "yes"

# Copy D4's value (keeps its quote-prefixed "1" text formatting) into D5
$ws3.Range("D4").Copy()
$ws3.Range("D5").PasteSpecial(-4163) # xlPasteValues

$ws3.Range("D5").Select()
